$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Headers for the new team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from an existing header cell
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Team record values for every data row (2 through 54)
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 68   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 94   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
